$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.911.21'
$ws.Range('E2').Value = '  +0.14%  '

$ws.Range('D3').Value = '2.572.59'
$ws.Range('E3').Value = '  +1.53%  '

$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').Value = '''313.13'
$ws.Range('E5').Value = '  -0.72%  '

$ws.Range('D6').Value = '''99.53'
$ws.Range('E6').Value = '  +3.52%  '

$ws.Range('E7').Value = '  -0.39%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('E9').Value = '  +0.24%  '

$ws.Range('D10').Value = '''35.96'
$ws.Range('E10').Value = '  -0.66%  '

$ws.Range('E11').Value = '  +0.44%  '

$ws.Range('D12').Value = '''7.47'
$ws.Range('E12').Value = '  -1.42%  '

$ws.Range('D13').Value = '2.967.09'
$ws.Range('E13').Value = '  +1.59%  '

$ws.Range('E14').Value = '  -1.27%  '

$ws.Range('D15').Value = '''15.96'
$ws.Range('E15').Value = '  +4.83%  '

$ws.Range('D16').Value = '2.572.20'
$ws.Range('E16').Value = '  +0.71%  '

$ws.Range('D17').Value = '''0.848'
$ws.Range('E17').Value = '  -0.50%  '

$ws.Range('D18').Value = '42.950.40'
$ws.Range('E18').Value = '  +0.14%  '

$ws.Range('E19').Value = '  -1.04%  '

$ws.Range('D20').Value = '''12.58'
$ws.Range('E20').Value = '  -4.03%  '

$ws.Range('E21').Value = '  -0.13%  '

$ws.Range('D22').Value = '''69.63'
$ws.Range('E22').Value = '  -0.57%  '

$ws.Range('D23').Value = '''250.28'
$ws.Range('E23').Value = '  -1.26%  '

$ws.Range('D24').Value = '''2.96'
$ws.Range('E24').Value = '  +0.52%  '

$ws.Range('E25').Value = '  -0.08%  '

$ws.Range('D26').Value = '''27.19'
$ws.Range('E26').Value = '  +2.00%  '

$ws.Range('E28').Value = '  -0.75%  '

$ws.Range('D29').Value = '''39.96'
$ws.Range('E29').Value = '  -1.52%  '

$ws.Range('D30').Value = '''10.29'
$ws.Range('E30').Value = '  -1.40%  '

$ws.Range('D31').Value = '''158.62'
$ws.Range('E31').Value = '  +0.40%  '

$ws.Range('E32').Value = '  -2.16%  '

$ws.Range('E33').Value = '  +1.00%  '

$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '''2.12'
$ws.Range('E34').Value = '  -2.17%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '''0.0803'
$ws.Range('E35').Value = '  +2.55%  '

$ws.Range('D36').Value = '''2.67'
$ws.Range('E36').Value = '  +0.09%  '

$ws.Range('D37').Value = '''18.70'
$ws.Range('E37').Value = '  -2.36%  '

$ws.Range('E38').Value = '  +11.24%  '

$ws.Range('E39').Value = '  +0.03%  '

$ws.Range('E40').Value = '  -0.21%  '

$ws.Range('D41').Value = '''23.56'
$ws.Range('E41').Value = '  +1.59%  '

$ws.Range('E42').Value = '  +7.29%  '

$ws.Range('E43').Value = '  -0.55%  '

$ws.Range('E44').Value = '  -0.07%  '

$ws.Range('E45').Value = '  -1.84%  '

$ws.Range('D46').Value = '1.998.94'
$ws.Range('E46').Value = '  -1.63%  '

$ws.Range('D47').Value = '''9.03'
$ws.Range('E47').Value = '  -1.91%  '

$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.817.52'
$ws.Range('E48').Value = '  +1.56%  '

$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '''0.198'
$ws.Range('E49').Value = '  +2.69%  '

$ws.Range('D50').Value = '''82.16'
$ws.Range('E50').Value = '  -3.51%  '

$ws.Range('D51').Value = '''74.77'
$ws.Range('E51').Value = '  -0.08%  '
